$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New localization entry for "ea45270a-3667-44e7-b006-483281ce8c47.md"
# Adds row 9 to all three sheets (Overview, zh-cn, de-de), expanding each
# sheet's table by one row.
# ---------------------------------------------------------------------------

# --- Overview sheet --------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A9").Value = "ea45270a-3667-44e7-b006-483281ce8c47.md"
$wsOverview.Range("B9").Value = "e2e\ea45270a-3667-44e7-b006-483281ce8c47.md"
$wsOverview.Range("C9").Value = ".md"
$wsOverview.Range("E9").Value = "Ready for handoff"
$wsOverview.Range("F9").Value = "Ready for handoff"
$wsOverview.Range("G9").Value = "2016-08-29 10:44:40"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/ea45270a-3667-44e7-b006-483281ce8c47.md",
    "",
    "",
    "e2e\ea45270a-3667-44e7-b006-483281ce8c47.md"
) | Out-Null

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A9").Value = "ea45270a-3667-44e7-b006-483281ce8c47.md"
$wsZhCn.Range("B9").Value = ".md"
$wsZhCn.Range("C9").Value = "Ready for handoff"
$wsZhCn.Range("D9").Value = "e2e"
$wsZhCn.Range("E9").Value = "ht"
$wsZhCn.Range("F9").Value = "False"
$wsZhCn.Range("G9").Value = "ea45270a-3667-44e7-b006-483281ce8c47.c85c050c1f71dfc97d5ac488a5247a799bac254a.zh-cn.xlf"
$wsZhCn.Range("H9").Value = "2016-08-29 10:44:35"
$wsZhCn.Range("K9").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M9").Value = "True"
$wsZhCn.Range("O9").Value = "False"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/ea45270a-3667-44e7-b006-483281ce8c47.md",
    "",
    "",
    "ea45270a-3667-44e7-b006-483281ce8c47.md"
) | Out-Null

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A9").Value = "ea45270a-3667-44e7-b006-483281ce8c47.md"
$wsDeDe.Range("B9").Value = ".md"
$wsDeDe.Range("C9").Value = "Ready for handoff"
$wsDeDe.Range("D9").Value = "e2e"
$wsDeDe.Range("E9").Value = "ht"
$wsDeDe.Range("F9").Value = "False"
$wsDeDe.Range("G9").Value = "ea45270a-3667-44e7-b006-483281ce8c47.c85c050c1f71dfc97d5ac488a5247a799bac254a.de-de.xlf"
$wsDeDe.Range("H9").Value = "2016-08-29 10:44:40"
$wsDeDe.Range("K9").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M9").Value = "True"
$wsDeDe.Range("O9").Value = "False"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/ea45270a-3667-44e7-b006-483281ce8c47.md",
    "",
    "",
    "ea45270a-3667-44e7-b006-483281ce8c47.md"
) | Out-Null
